$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking text values in D/E columns from Excel auto-number-conversion
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "62.777.86"
$ws.Range("D3").Value = "3.469.93"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "414.72"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("D6").Value = "130.92"
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").Value = "42.72"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "9.73"
$ws.Range("E12").Value = "  +6.68%  "
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").Value = "4.014.69"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").Value = "0.141"
$ws.Range("D16").Value = "20.54"
$ws.Range("E16").Value = "  -4.08%  "
$ws.Range("D17").Value = "3.463.05"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "62.715.61"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "469.34"
$ws.Range("E21").Value = "  +3.95%  "
$ws.Range("D22").Value = "90.89"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("D23").Value = "3.27"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").Value = "13.36"
$ws.Range("E24").Value = "  +2.99%  "
$ws.Range("D25").Value = "10.62"
$ws.Range("E25").Value = "  +21.25%  "
$ws.Range("D26").Value = "3.32"
$ws.Range("D27").Value = "33.39"
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("D28").Value = "4.81"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "2.66"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("D34").Value = "41.14"
$ws.Range("E34").Value = "  -4.48%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "58.76"
$ws.Range("E36").Value = "  +8.42%  "
$ws.Range("D37").Value = "0.0490"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "3.08"
$ws.Range("E38").Value = "  +5.10%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "3.35"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("E43").Value = "  +6.78%  "
$ws.Range("D44").Value = "145.93"
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D45").Value = "4.37"
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("D46").Value = "2.08"
$ws.Range("E46").Value = "  +4.26%  "
$ws.Range("E47").Value = "  +11.28%  "
$ws.Range("D48").Value = "0.0₃0565"
$ws.Range("E48").Value = "  +34.84%  "
$ws.Range("D49").Value = "16.40"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("D50").Value = "22.36"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  -0.31%  "

# Restore original (default) cell style now that text values are set
$ws.Range("D2:E51").Style = "Normal"
